# NIT-9004205577.xlsx - "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# Adds two new worker rows (MAYERLI BARRERA BARRERA, doc 1047488141) to the
# "Estado de Cuenta" table, swaps the two existing "Periodo Mora" values for
# the first worker (YAJAIRA LUZ BOTELLO RINCONES), and updates the summary
# totals (Valor Mora, Cant. Trabajadores, Cant. Periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room: insert two blank rows right before the current last
#    data row (18) - this pushes the old footer rows (22/23 -> 24/25)
#    down automatically and keeps their merged cells intact.
# ------------------------------------------------------------------
$ws.Range("B18:J19").Insert()

# ------------------------------------------------------------------
# 2) Fix up formatting for the table rows:
#    - row 17 currently still carries the "last row" (thicker border)
#      style; row 16's "middle row" style needs to be re-applied to it
#      (and to the newly inserted row 18).
#    - the newly inserted row 19 needs the "last row" style that used
#      to belong to row 17.
# ------------------------------------------------------------------
$ws.Range("B17:J17").Copy($ws.Range("B19:J19"))
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# ------------------------------------------------------------------
# 3) Data: the two existing "Periodo Mora" values (1709/1710) swap
#    between row 16 and row 17.
# ------------------------------------------------------------------
$ws.Range("E16").Value = "1710"
$ws.Range("E17").Value = "1709"

# ------------------------------------------------------------------
# 4) Data: two brand-new rows for worker MAYERLI BARRERA BARRERA.
# ------------------------------------------------------------------
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047488141"
$ws.Range("D18").Value = "MAYERLI BARRERA BARRERA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 60000
$ws.Range("G18").Value = 1500000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047488141"
$ws.Range("D19").Value = "MAYERLI BARRERA BARRERA"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 60000
$ws.Range("G19").Value = 1500000

# ------------------------------------------------------------------
# 5) Summary block updates.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 192000   # VALOR MORA total
$ws.Range("C13").Value = 2        # Cant. Trabajadores
$ws.Range("F13").Value = 4        # Cant. Periodos

Write-Output "edit applied"
